# LST 1469 upload budgets from adjustment columns
# Add three new "adjustment" header columns (ADJ1, ADJ2, ADJ3) to the
# Budgets sheet, right after the existing "Project" / month columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budgets")

# New header cells: S1=ADJ1, T1=ADJ2, U1=ADJ3
$ws.Cells.Item(1, 19).Value2 = "ADJ1"
$ws.Cells.Item(1, 20).Value2 = "ADJ2"
$ws.Cells.Item(1, 21).Value2 = "ADJ3"

# Reflect the selection/active cell state left behind after adding the
# columns (user ended up with S4 selected).
$ws.Range("S4").Select() | Out-Null
